$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0476088747382164
$ws.Range("C2").Value = 0.013125336728990078
$ws.Range("D2").Value = 0.007366156205534935
$ws.Range("E2").Value = 0.004317492712289095
$ws.Range("F2").Value = 0.00000005809338432527511
$ws.Range("G2").Value = 0.001304448815062642
$ws.Range("J2").Value = 0.1276422142982483
$ws.Range("K2").Value = 1.4611644744873047
